# Apply the data-entry changes to the "jobs" workbook.
#
# Only four plain input cells actually change value; every other cell that
# differs in the target diff (A39, B40, E40, F40, G40, B42, C42, E42, F42,
# G42, B44, C44, E44, B46, C46, E46, C49) is a formula that recalculates
# automatically once these inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A40").Value = 2
$ws.Range("A42").Value = 15
$ws.Range("C45").Value = 0.5
$ws.Range("E45").Value = 0.2

# Restore the on-screen scroll/selection state captured in the target file.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = $ws.Range("B30").Row
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("B30").Column
$ws.Range("C45").Select()
